# Auto-generated: restore correct per-row betting-odds data for the
# Estonia Meistriliiga sheet. Several rows' match records were reordered
# (the match `id`, teams, and odds columns B:AC move together as a unit;
# only column A, the running row/match counter, stays put).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$row7 = New-Object 'object[,]' 1,28
$row7[0,0] = 6397040
$row7[0,1] = 'Estonia Meistriliiga'
$row7[0,2] = 'Estonia Meistriliiga'
$row7[0,3] = 45000.54166666666
$row7[0,4] = 'Harju JK Laagri'
$row7[0,5] = 'JK Tallinna Kalev'
$row7[0,6] = 0
$row7[0,7] = 2
$row7[0,8] = 'A'
$row7[0,9] = 2.2
$row7[0,10] = 3.4
$row7[0,11] = 2.75
$row7[0,12] = 2.25
$row7[0,13] = 3.4
$row7[0,14] = 2.625
$row7[0,15] = 0
$row7[0,16] = 1.8
$row7[0,17] = 2
$row7[0,18] = 2.75
$row7[0,19] = 1.925
$row7[0,20] = 1.875
$row7[0,21] = -1
$row7[0,22] = -1
$row7[0,23] = 1.625
$row7[0,24] = -1
$row7[0,25] = 1
$row7[0,26] = -1
$row7[0,27] = 0.875
$ws.Range("B7:AC7").Value2 = $row7

# Row 8
$row8 = New-Object 'object[,]' 1,28
$row8[0,0] = 6394634
$row8[0,1] = 'Estonia Meistriliiga'
$row8[0,2] = 'Estonia Meistriliiga'
$row8[0,3] = 45000.54166666666
$row8[0,4] = 'FC Levadia Tallinn'
$row8[0,5] = 'JK Tammeka Tartu'
$row8[0,6] = 3
$row8[0,7] = 0
$row8[0,8] = 'H'
$row8[0,9] = 1.333
$row8[0,10] = 4.5
$row8[0,11] = 7
$row8[0,12] = 1.4
$row8[0,13] = 4.2
$row8[0,14] = 5.75
$row8[0,15] = -1.25
$row8[0,16] = 2
$row8[0,17] = 1.8
$row8[0,18] = 2.75
$row8[0,19] = 1.95
$row8[0,20] = 1.85
$row8[0,21] = 0.3999999999999999
$row8[0,22] = -1
$row8[0,23] = -1
$row8[0,24] = 1
$row8[0,25] = -1
$row8[0,26] = 0.475
$row8[0,27] = -0.5
$ws.Range("B8:AC8").Value2 = $row8

# Row 77
$row77 = New-Object 'object[,]' 1,28
$row77[0,0] = 6139018
$row77[0,1] = 'Estonia Meistriliiga'
$row77[0,2] = 'Estonia Meistriliiga'
$row77[0,3] = 45084.5
$row77[0,4] = 'JK Tallinna Kalev'
$row77[0,5] = 'JK Trans Narva'
$row77[0,6] = 0
$row77[0,7] = 1
$row77[0,8] = 'A'
$row77[0,9] = 2.4
$row77[0,10] = 3.4
$row77[0,11] = 2.5
$row77[0,12] = 2.875
$row77[0,13] = 3.1
$row77[0,14] = 2.3
$row77[0,15] = 0.25
$row77[0,16] = 1.75
$row77[0,17] = 2.05
$row77[0,18] = 2.25
$row77[0,19] = 1.925
$row77[0,20] = 1.875
$row77[0,21] = -1
$row77[0,22] = -1
$row77[0,23] = 1.3
$row77[0,24] = -1
$row77[0,25] = 1.05
$row77[0,26] = -1
$row77[0,27] = 0.875
$ws.Range("B77:AC77").Value2 = $row77

# Row 78
$row78 = New-Object 'object[,]' 1,28
$row78[0,0] = 6139017
$row78[0,1] = 'Estonia Meistriliiga'
$row78[0,2] = 'Estonia Meistriliiga'
$row78[0,3] = 45084.5
$row78[0,4] = 'JK Tammeka Tartu'
$row78[0,5] = 'Harju JK Laagri'
$row78[0,6] = 2
$row78[0,7] = 0
$row78[0,8] = 'H'
$row78[0,9] = 1.666
$row78[0,10] = 3.6
$row78[0,11] = 4.2
$row78[0,12] = 1.727
$row78[0,13] = 3.5
$row78[0,14] = 4
$row78[0,15] = -0.75
$row78[0,16] = 2
$row78[0,17] = 1.8
$row78[0,18] = 2.5
$row78[0,19] = 1.9
$row78[0,20] = 1.9
$row78[0,21] = 0.7270000000000001
$row78[0,22] = -1
$row78[0,23] = -1
$row78[0,24] = 1
$row78[0,25] = -1
$row78[0,26] = -1
$row78[0,27] = 0.8999999999999999
$ws.Range("B78:AC78").Value2 = $row78

# Row 168
$row168 = New-Object 'object[,]' 1,28
$row168[0,0] = 6416370
$row168[0,1] = 'Estonia Meistriliiga'
$row168[0,2] = 'Estonia Meistriliiga'
$row168[0,3] = 45231.54166666666
$row168[0,4] = 'FC Levadia Tallinn'
$row168[0,5] = 'Parnu JK Vaprus'
$row168[0,6] = 0
$row168[0,7] = 0
$row168[0,8] = 'D'
$row168[0,9] = 1.166
$row168[0,10] = 7
$row168[0,11] = 11
$row168[0,12] = 1.2
$row168[0,13] = 6
$row168[0,14] = 11
$row168[0,15] = -2
$row168[0,16] = 1.85
$row168[0,17] = 1.95
$row168[0,18] = 3
$row168[0,19] = 1.85
$row168[0,20] = 1.95
$row168[0,21] = -1
$row168[0,22] = 5
$row168[0,23] = -1
$row168[0,24] = -1
$row168[0,25] = 0.95
$row168[0,26] = -1
$row168[0,27] = 0.95
$ws.Range("B168:AC168").Value2 = $row168

# Row 169
$row169 = New-Object 'object[,]' 1,28
$row169[0,0] = 6482819
$row169[0,1] = 'Estonia Meistriliiga'
$row169[0,2] = 'Estonia Meistriliiga'
$row169[0,3] = 45231.54166666666
$row169[0,4] = 'JK Tammeka Tartu'
$row169[0,5] = 'FC Kuressaare'
$row169[0,6] = 0
$row169[0,7] = 1
$row169[0,8] = 'A'
$row169[0,9] = 1.833
$row169[0,10] = 3.5
$row169[0,11] = 3.5
$row169[0,12] = 2.1
$row169[0,13] = 3.4
$row169[0,14] = 2.875
$row169[0,15] = -0.25
$row169[0,16] = 1.975
$row169[0,17] = 1.825
$row169[0,18] = 3
$row169[0,19] = 1.825
$row169[0,20] = 1.975
$row169[0,21] = -1
$row169[0,22] = -1
$row169[0,23] = 1.875
$row169[0,24] = -1
$row169[0,25] = 0.825
$row169[0,26] = -1
$row169[0,27] = 0.9750000000000001
$ws.Range("B169:AC169").Value2 = $row169

# Row 177
$row177 = New-Object 'object[,]' 1,28
$row177[0,0] = 6533597
$row177[0,1] = 'Estonia Meistriliiga'
$row177[0,2] = 'Estonia Meistriliiga'
$row177[0,3] = 45241.375
$row177[0,4] = 'FC Kuressaare'
$row177[0,5] = 'Parnu JK Vaprus'
$row177[0,6] = 1
$row177[0,7] = 0
$row177[0,8] = 'H'
$row177[0,9] = 2.5
$row177[0,10] = 3.4
$row177[0,11] = 2.5
$row177[0,12] = 2.15
$row177[0,13] = 3.6
$row177[0,14] = 2.875
$row177[0,15] = -0.25
$row177[0,16] = 1.95
$row177[0,17] = 1.85
$row177[0,18] = 2.75
$row177[0,19] = 1.95
$row177[0,20] = 1.85
$row177[0,21] = 1.15
$row177[0,22] = -1
$row177[0,23] = -1
$row177[0,24] = 0.95
$row177[0,25] = -1
$row177[0,26] = -1
$row177[0,27] = 0.8500000000000001
$ws.Range("B177:AC177").Value2 = $row177

# Row 178
$row178 = New-Object 'object[,]' 1,28
$row178[0,0] = 6537957
$row178[0,1] = 'Estonia Meistriliiga'
$row178[0,2] = 'Estonia Meistriliiga'
$row178[0,3] = 45241.375
$row178[0,4] = 'FC Flora Tallinn'
$row178[0,5] = 'JK Nomme Kalju'
$row178[0,6] = 0
$row178[0,7] = 0
$row178[0,8] = 'D'
$row178[0,9] = 1.4
$row178[0,10] = 4
$row178[0,11] = 7.5
$row178[0,12] = 1.5
$row178[0,13] = 4.2
$row178[0,14] = 5
$row178[0,15] = -1
$row178[0,16] = 1.85
$row178[0,17] = 1.95
$row178[0,18] = 2.75
$row178[0,19] = 1.85
$row178[0,20] = 1.95
$row178[0,21] = -1
$row178[0,22] = 3.2
$row178[0,23] = -1
$row178[0,24] = -1
$row178[0,25] = 0.95
$row178[0,26] = -1
$row178[0,27] = 0.95
$ws.Range("B178:AC178").Value2 = $row178

# Row 179
$row179 = New-Object 'object[,]' 1,28
$row179[0,0] = 6537869
$row179[0,1] = 'Estonia Meistriliiga'
$row179[0,2] = 'Estonia Meistriliiga'
$row179[0,3] = 45241.375
$row179[0,4] = 'JK Tallinna Kalev'
$row179[0,5] = 'JK Trans Narva'
$row179[0,6] = 5
$row179[0,7] = 0
$row179[0,8] = 'H'
$row179[0,9] = 1.6
$row179[0,10] = 4
$row179[0,11] = 4.5
$row179[0,12] = 1.65
$row179[0,13] = 4
$row179[0,14] = 4.333
$row179[0,15] = -0.75
$row179[0,16] = 1.8
$row179[0,17] = 2
$row179[0,18] = 2.75
$row179[0,19] = 1.9
$row179[0,20] = 1.9
$row179[0,21] = 0.6499999999999999
$row179[0,22] = -1
$row179[0,23] = -1
$row179[0,24] = 0.8
$row179[0,25] = -1
$row179[0,26] = 0.8999999999999999
$row179[0,27] = -1
$ws.Range("B179:AC179").Value2 = $row179

# Row 180
$row180 = New-Object 'object[,]' 1,28
$row180[0,0] = 6535416
$row180[0,1] = 'Estonia Meistriliiga'
$row180[0,2] = 'Estonia Meistriliiga'
$row180[0,3] = 45241.375
$row180[0,4] = 'Paide Linnameeskond'
$row180[0,5] = 'FC Levadia Tallinn'
$row180[0,6] = 2
$row180[0,7] = 2
$row180[0,8] = 'D'
$row180[0,9] = 3
$row180[0,10] = 3.8
$row180[0,11] = 2
$row180[0,12] = 3
$row180[0,13] = 4
$row180[0,14] = 1.909
$row180[0,15] = 0.5
$row180[0,16] = 1.85
$row180[0,17] = 1.95
$row180[0,18] = 2.75
$row180[0,19] = 1.95
$row180[0,20] = 1.85
$row180[0,21] = -1
$row180[0,22] = 3
$row180[0,23] = -1
$row180[0,24] = 0.8500000000000001
$row180[0,25] = -1
$row180[0,26] = 0.95
$row180[0,27] = -1
$ws.Range("B180:AC180").Value2 = $row180

